$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the per-row "pass/fail" and "expected result" helper columns (F, G)
# for every data row, but keep the header labels in F1/G1.
$ws.Range("F2:G18").ClearContents() | Out-Null

# Update the selection shown when the sheet was last saved
$ws.Range("I14").Select() | Out-Null
